$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '34.614.42'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  +1.27%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.796.95'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  +0.86%  '
$ws.Range('E4').Value = '  -0.09%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '227.03'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.42%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.560'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +2.40%  '
$ws.Range('E7').Value = '  -0.13%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '33.01'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +3.96%  '
$ws.Range('E9').Value = '  +1.81%  '
$ws.Range('E10').Value = '  +0.83%  '
$ws.Range('E11').Value = '  +0.40%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '2.055.41'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +0.80%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '1.827.61'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +2.66%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '11.07'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +0.64%  '
$ws.Range('E15').Value = '  +2.51%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '34.571.69'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +1.29%  '
$ws.Range('E17').Value = '  +2.95%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '68.92'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +1.51%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '247.53'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +0.48%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '11.32'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +3.01%  '
$ws.Range('E22').Value = '  -0.14%  '
$ws.Range('E23').Value = '  +2.02%  '
$ws.Range('E24').Value = '  +1.55%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '166.51'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +2.47%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '7.31'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +1.83%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '16.60'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +1.88%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '0.117'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +2.43%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.00'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -0.04%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '4.08'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +11.27%  '
$ws.Range('B31').NumberFormat = '@'
$ws.Range('B31').Value = 'PancakeSwap'
$ws.Range('B31').Style = 'Normal'
$ws.Range('C31').NumberFormat = '@'
$ws.Range('C31').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('C31').Style = 'Normal'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.24'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +1.11%  '
$ws.Range('B32').NumberFormat = '@'
$ws.Range('B32').Value = 'Hedera'
$ws.Range('B32').Style = 'Normal'
$ws.Range('C32').NumberFormat = '@'
$ws.Range('C32').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('C32').Style = 'Normal'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.0525'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +1.07%  '
$ws.Range('B33').NumberFormat = '@'
$ws.Range('B33').Value = 'Filecoin'
$ws.Range('B33').Style = 'Normal'
$ws.Range('C33').NumberFormat = '@'
$ws.Range('C33').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('C33').Style = 'Normal'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '3.81'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +2.28%  '
$ws.Range('E34').Value = '  +2.85%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.428.61'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -0.94%  '
$ws.Range('E36').Value = '  +7.99%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.672'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +2.51%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '1.06'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +2.03%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.0192'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +0.83%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '85.89'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +6.95%  '
$ws.Range('E41').Value = '  +1.68%  '
$ws.Range('B42').NumberFormat = '@'
$ws.Range('B42').Value = 'MXToken'
$ws.Range('B42').Style = 'Normal'
$ws.Range('C42').NumberFormat = '@'
$ws.Range('C42').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('C42').Style = 'Normal'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '2.76'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +3.13%  '
$ws.Range('B43').NumberFormat = '@'
$ws.Range('B43').Value = 'ARBITRUM'
$ws.Range('B43').Style = 'Normal'
$ws.Range('C43').NumberFormat = '@'
$ws.Range('C43').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('C43').Style = 'Normal'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.935'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +1.13%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '13.73'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +1.13%  '
$ws.Range('E45').Value = '  +3.89%  '
$ws.Range('E46').Value = '  +0.95%  '
$ws.Range('E47').Value = '  +0.58%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.955.25'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +0.73%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '106.08'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +1.34%  '
$ws.Range('E50').Value = '  -0.11%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.0₆0128'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -5.54%  '
